$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '66.352.61'
$ws.Range("E2").Value = '  +6.24%  '

$ws.Range("D3").Value = '3.549.30'
$ws.Range("E3").Value = '  +3.10%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = '  -0.10%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '418.16'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +1.08%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '130.69'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  +0.77%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.656'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  +5.06%  '

$ws.Range("D8").Value = '3.545.04'
$ws.Range("E8").Value = '  +3.21%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.999'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  -0.07%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.777'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  +7.10%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.175'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  +24.53%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0000308'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  +40.83%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '43.13'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  +1.21%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '10.02'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  +5.87%  '

$ws.Range("D15").Value = '4.115.70'
$ws.Range("E15").Value = '  +3.24%  '

$ws.Range("E16").Value = '  -0.15%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '20.49'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  +0.11%  '

$ws.Range("D18").Value = '3.584.15'
$ws.Range("E18").Value = '  +3.74%  '

$ws.Range("E19").Value = '  +4.83%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '12.50'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  -4.67%  '

$ws.Range("D21").Value = '66.272.57'
$ws.Range("E21").Value = '  +6.09%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '453.61'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  -4.80%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '90.24'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  -0.85%  '

$ws.Range("E24").Value = '  -2.23%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '13.17'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  -1.66%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.37'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  +1.96%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.04'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  -4.58%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '34.68'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  +4.01%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '4.84'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  +0.70%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '12.48'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  +4.45%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.79'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  +5.25%  '

$ws.Range("E32").Value = '  +5.67%  '

$ws.Range("E33").Value = '  -3.63%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.160'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  -3.74%  '

$ws.Range("E35").Value = '  -0.31%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '39.09'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  -3.79%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '57.02'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  -2.80%  '

$ws.Range("E38").Value = '  +42.42%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0500'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  +2.10%  '

$ws.Range("E40").Value = '  +10.23%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.998'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  -0.12%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.77'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  +3.18%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '3.01'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -0.25%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '149.04'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  +2.82%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '4.40'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  +1.30%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.25'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  -2.79%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.309'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -4.51%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.99'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  -3.79%  '

$ws.Range("E49").Value = '  -3.67%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.144'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  +4.55%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '15.53'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -5.12%  '
